$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row for Albahaca (Agrícola del Norte S.A. de Arica)
# is inserted at row 5, pushing the existing rows 5-32 down to rows 6-33.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new week's data. The
# non-varying descriptive columns match every other row in this sheet.
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(5, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(5, 4).Value = 44635
$ws.Cells.Item(5, 5).Value = 15
$ws.Cells.Item(5, 6).Value = 100112052
$ws.Cells.Item(5, 7).Value = "Albahaca"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 300
$ws.Cells.Item(5, 11).Value = 1900
$ws.Cells.Item(5, 12).Value = 2000
$ws.Cells.Item(5, 13).Value = 1950
$ws.Cells.Item(5, 14).Value = "$/paquete"
$ws.Cells.Item(5, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(5, 16).Value = 1950
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = "Hortaliza"
